$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header E1 ("EQP" -> "Flag3") and add new header F1 ("Flag4")
$ws.Range("E1").Value = "Flag3"
$ws.Range("F1").Value = "Flag4"

# 2. Normalize cell styles: make every data cell B2:F7 use the centered style
#    (same alignment the header row / column A already use), which collapses
#    the old "style 2" xf record into the shared "style 1" xf record.
$ws.Range("B2:F7").HorizontalAlignment = -4108

# 3. Make sure the new column F header also carries the same style as the
#    rest of row 1.
$ws.Range("F1").HorizontalAlignment = -4108

# 4. Row heights for rows 1-5 recompute to 12.8 (matching rows 6-7 already
#    present in the original file).
$ws.Rows.Item(1).RowHeight = 12.8
$ws.Rows.Item(2).RowHeight = 12.8
$ws.Rows.Item(3).RowHeight = 12.8
$ws.Rows.Item(4).RowHeight = 12.8
$ws.Rows.Item(5).RowHeight = 12.8

# 5. Update the active selection to C2.
$ws.Range("C2").Select()

# 6. Data validation: the "list" validation used to cover the header row
#    (A1:IQ1) and column A (A2:A1007). After adding the new column the two
#    new header cells (E1:F1) no longer carry the validation, while the new
#    cells D2:F2 and E3:F7 gain it.
$ws.Range("E1:F1").Validation.Delete()

$rngD2F2 = $ws.Range("D2:F2")
$rngD2F2.Validation.Add(3, 1, 1, '"qwer%yuiop_1234567890"', 0)
$dv1 = $rngD2F2.Validation
$dv1.ErrorTitle = "Not Applicable"
$dv1.ErrorMessage = "Cannot change the value"
$dv1.IgnoreBlank = $true
$dv1.InCellDropdown = $false
$dv1.ShowError = $true
$dv1.ShowInput = $false

$rngE3F7 = $ws.Range("E3:F7")
$rngE3F7.Validation.Add(3, 1, 1, '"qwer%yuiop_1234567890"', 0)
$dv2 = $rngE3F7.Validation
$dv2.ErrorTitle = "Not Applicable"
$dv2.ErrorMessage = "Cannot change the value"
$dv2.IgnoreBlank = $true
$dv2.InCellDropdown = $false
$dv2.ShowError = $true
$dv2.ShowInput = $false
